$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 62.444443
$ws.Range("I11").Value = 62.444443
$ws.Range("K11").Value = 62.444443
$ws.Range("M11").Value = 77.55555699999999
$ws.Range("H33").Value = 1151.7059
$ws.Range("I33").Value = 1477.6666
$ws.Range("J33").Value = 785
$ws.Range("K33").Value = 1477.6666
$ws.Range("L33").Value = 785
$ws.Range("M33").Value = -1248.6666
$ws.Range("N33").Value = -1243
$ws.Range("H64").Value = 3835.7144
$ws.Range("I64").Value = 3700
$ws.Range("K64").Value = 3700
$ws.Range("M64").Value = -3452
$ws.Range("H67").Value = 3835.7144
$ws.Range("I67").Value = 3700
$ws.Range("K67").Value = 3700
$ws.Range("M67").Value = -2842
$ws.Range("H74").Value = 4277.6665
$ws.Range("I74").Value = 4500.5
$ws.Range("J74").Value = 4214
$ws.Range("K74").Value = 4500.5
$ws.Range("L74").Value = 4214
$ws.Range("M74").Value = -3564.5
$ws.Range("N74").Value = -6086
$ws.Range("H76").Value = 3383.3333
$ws.Range("I76").Value = 3475
$ws.Range("K76").Value = 3475
$ws.Range("M76").Value = -3160
$ws.Range("H77").Value = 4277.6665
$ws.Range("I77").Value = 4500.5
$ws.Range("J77").Value = 4214
$ws.Range("K77").Value = 22502.5
$ws.Range("L77").Value = 21070
$ws.Range("M77").Value = -17822.5
$ws.Range("N77").Value = -30430
$ws.Range("H79").Value = 3383.3333
$ws.Range("I79").Value = 3475
$ws.Range("K79").Value = 3475
$ws.Range("M79").Value = -2383

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2611.8
$ws.Range("I2").Value = 2859.5557
$ws.Range("J2").Value = 2409.0908
$ws.Range("K2").Value = 2859.5557
$ws.Range("L2").Value = 2409.0908
$ws.Range("M2").Value = -2746.5557
$ws.Range("N2").Value = -2635.0908
$ws.Range("H45").Value = 1594.2354
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 1721.5714
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 1721.5714
$ws.Range("M45").Value = -623
$ws.Range("N45").Value = -2475.5714
$ws.Range("H63").Value = 3319.087
$ws.Range("I63").Value = 2783.6875
$ws.Range("J63").Value = 4542.857
$ws.Range("K63").Value = 2783.6875
$ws.Range("L63").Value = 4542.857
$ws.Range("M63").Value = -2097.6875
$ws.Range("N63").Value = -5914.857
$ws.Range("H66").Value = 3319.087
$ws.Range("I66").Value = 2783.6875
$ws.Range("J66").Value = 4542.857
$ws.Range("K66").Value = 13918.4375
$ws.Range("L66").Value = 22714.285
$ws.Range("M66").Value = -10486.4375
$ws.Range("N66").Value = -29578.285
$ws.Range("H88").Value = 2341.7273
$ws.Range("I88").Value = 2177.25
$ws.Range("J88").Value = 2435.7144
$ws.Range("K88").Value = 2177.25
$ws.Range("L88").Value = 2435.7144
$ws.Range("M88").Value = -1771.25
$ws.Range("N88").Value = -3247.7144
$ws.Range("H91").Value = 2341.7273
$ws.Range("I91").Value = 2177.25
$ws.Range("J91").Value = 2435.7144
$ws.Range("K91").Value = 2177.25
$ws.Range("L91").Value = 2435.7144
$ws.Range("M91").Value = -773.25
$ws.Range("N91").Value = -5243.7144
$ws.Range("H97").Value = 956.35
$ws.Range("I97").Value = 928.17645
$ws.Range("K97").Value = 928.17645
$ws.Range("M97").Value = -432.17645
$ws.Range("H110").Value = 1745.4
$ws.Range("I110").Value = 1425.2307
$ws.Range("J110").Value = 2340
$ws.Range("K110").Value = 1425.2307
$ws.Range("L110").Value = 2340
$ws.Range("M110").Value = 619.7692999999999
$ws.Range("N110").Value = -6430
$ws.Range("H116").Value = 2611.8
$ws.Range("I116").Value = 2859.5557
$ws.Range("J116").Value = 2409.0908
$ws.Range("K116").Value = 2859.5557
$ws.Range("L116").Value = 2409.0908
$ws.Range("M116").Value = -565.5556999999999
$ws.Range("N116").Value = -6997.0908

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2611.8
$ws.Range("I3").Value = 2859.5557
$ws.Range("J3").Value = 2409.0908
$ws.Range("K3").Value = 2859.5557
$ws.Range("L3").Value = 2409.0908
$ws.Range("M3").Value = -2745.5557
$ws.Range("N3").Value = -2637.0908
$ws.Range("H86").Value = 1500.4445
$ws.Range("I86").Value = 1426.579
$ws.Range("J86").Value = 1675.875
$ws.Range("K86").Value = 1426.579
$ws.Range("L86").Value = 1675.875
$ws.Range("M86").Value = -303.579
$ws.Range("N86").Value = -3921.875
$ws.Range("H89").Value = 1500.4445
$ws.Range("I89").Value = 1426.579
$ws.Range("J89").Value = 1675.875
$ws.Range("K89").Value = 7132.895
$ws.Range("L89").Value = 8379.375
$ws.Range("M89").Value = -1516.895
$ws.Range("N89").Value = -19611.375
$ws.Range("H94").Value = 767.4400000000001
$ws.Range("I94").Value = 707.2105
$ws.Range("J94").Value = 958.1667
$ws.Range("K94").Value = 707.2105
$ws.Range("L94").Value = 958.1667
$ws.Range("M94").Value = -256.2105
$ws.Range("N94").Value = -1860.1667
$ws.Range("H99").Value = 1774.7391
$ws.Range("I99").Value = 1352.7858
$ws.Range("K99").Value = 1352.7858
$ws.Range("M99").Value = 145.2141999999999
$ws.Range("H107").Value = 1547.3846
$ws.Range("I107").Value = 1501.8572
$ws.Range("J107").Value = 1600.5
$ws.Range("K107").Value = 1501.8572
$ws.Range("L107").Value = 1600.5
$ws.Range("M107").Value = 418.1428000000001
$ws.Range("N107").Value = -5440.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.666664
$ws.Range("I7").Value = 34
$ws.Range("K7").Value = 34
$ws.Range("M7").Value = 79
$ws.Range("H31").Value = 45455870
$ws.Range("I31").Value = 41667720
$ws.Range("J31").Value = 55557612
$ws.Range("K31").Value = 41667720
$ws.Range("L31").Value = 55557612
$ws.Range("M31").Value = -41667425
$ws.Range("N31").Value = -55558202
$ws.Range("H34").Value = 45455870
$ws.Range("I34").Value = 41667720
$ws.Range("J34").Value = 55557612
$ws.Range("K34").Value = 41667720
$ws.Range("L34").Value = 55557612
$ws.Range("M34").Value = -41667518
$ws.Range("N34").Value = -55558016
$ws.Range("H62").Value = 3045.4546
$ws.Range("I62").Value = 2312.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2312.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1688.5
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3045.4546
$ws.Range("I65").Value = 2312.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 11562.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -8442.5
$ws.Range("N65").Value = -31240
$ws.Range("H105").Value = 1050
$ws.Range("I105").Value = 1050
$ws.Range("K105").Value = 1050
$ws.Range("M105").Value = 697
$ws.Range("H122").Value = 1726.4
$ws.Range("I122").Value = 1430
$ws.Range("J122").Value = 2022.8
$ws.Range("K122").Value = 4290
$ws.Range("L122").Value = 6068.4
$ws.Range("M122").Value = -1840
$ws.Range("N122").Value = -10968.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 567.375
$ws.Range("I113").Value = 591.5909
$ws.Range("J113").Value = 537.7778
$ws.Range("K113").Value = 1774.7727
$ws.Range("L113").Value = 1613.3334
$ws.Range("M113").Value = 395.2273
$ws.Range("N113").Value = -5953.3334
$ws.Range("H132").Value = 3196.2666
$ws.Range("I132").Value = 2994.9092
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 26954.1828
$ws.Range("L132").Value = 33750
$ws.Range("M132").Value = -24424.1828
$ws.Range("N132").Value = -38810

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5228.6514
$ws.Range("I70").Value = 5079.143
$ws.Range("J70").Value = 5507.7334
$ws.Range("K70").Value = 5079.143
$ws.Range("L70").Value = 5507.7334
$ws.Range("M70").Value = -4809.143
$ws.Range("N70").Value = -6047.7334
$ws.Range("H73").Value = 5228.6514
$ws.Range("I73").Value = 5079.143
$ws.Range("J73").Value = 5507.7334
$ws.Range("K73").Value = 5079.143
$ws.Range("L73").Value = 5507.7334
$ws.Range("M73").Value = -4143.143
$ws.Range("N73").Value = -7379.7334
$ws.Range("H80").Value = 2766.3
$ws.Range("J80").Value = 3431.5
$ws.Range("L80").Value = 3431.5
$ws.Range("N80").Value = -5427.5
$ws.Range("H83").Value = 2766.3
$ws.Range("J83").Value = 3431.5
$ws.Range("L83").Value = 17157.5
$ws.Range("N83").Value = -27141.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 50005
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31996
$ws.Range("H75").Value = 29000
$ws.Range("J75").Value = 29000
$ws.Range("L75").Value = 29000
$ws.Range("N75").Value = -30872
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99984
$ws.Range("H78").Value = 29000
$ws.Range("J78").Value = 29000
$ws.Range("L78").Value = 87000
$ws.Range("N78").Value = -96360
$ws.Range("H93").Value = 26461.857
$ws.Range("I93").Value = 1614.5217
$ws.Range("J93").Value = 140759.6
$ws.Range("K93").Value = 1614.5217
$ws.Range("L93").Value = 140759.6
$ws.Range("M93").Value = -366.5217
$ws.Range("N93").Value = -143255.6
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180
$ws.Range("H122").Value = 4823.25
$ws.Range("I122").Value = 1493
$ws.Range("J122").Value = 5933.3335
$ws.Range("K122").Value = 4479
$ws.Range("L122").Value = 17800.0005
$ws.Range("M122").Value = -2029
$ws.Range("N122").Value = -22700.0005
